$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the Category Code value on row 2 (was "null#vital-signs")
$ws.Range("C2").Value = "Observation Category Codes#vital-signs"

# Append a new profile row (US Core Laboratory Result Observation Profile)
$ws.Range("A5").Value = "us-core-observation-lab"
$ws.Range("B5").Value = "US Core Laboratory Result Observation Profile"
$ws.Range("C5").Value = "Observation Category Codes#laboratory"
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = "http://hl7.org/fhir/us/core/ValueSet/us-core-laboratory-test-codes (extensible)"
$ws.Range("G5").Value = "dateTimeĵ, Periodĵ, Timingĵ, instantĵ"
$ws.Range("H5").Value = "Quantityĵ, CodeableConceptĵ, stringĵ, booleanĵ, integerĵ, Rangeĵ, Ratioĵ, SampledDataĵ, timeĵ, dateTimeĵ, Periodĵ"
$ws.Range("I5").Value = "optional"
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = ""

# Match the formatting used by the other data rows (copy format only, xlPasteFormats = -4122)
$ws.Range("A4:K4").Copy()
$ws.Range("A5:K5").PasteSpecial(-4122)
